{"js": "// Replace the math expressions in the worksheet table, cell by cell,\n// preserving each run's existing formatting (font, size, paragraph\n// justification, etc). Each entry is [rowIndex, colIndex, oldText, newText]\n// (0-based row/col), taken from the authoritative edit.\nconst edits = [\n  [0, 0, \"75+19=\", \"66+25=\"],\n  [0, 1, \"99-43=\", \"72-6=\"],\n  [0, 2, \"89-3=\", \"2+97=\"],\n  [0, 3, \"12+66=\", \"79+17=\"],\n  [0, 4, \"90-6=\", \"7+67=\"],\n  [1, 0, \"4+25=\", \"93-24=\"],\n  [1, 1, \"89-78=\", \"39+29=\"],\n  [1, 2, \"74-25=\", \"82-1=\"],\n  [1, 3, \"53-34=\", \"46+51=\"],\n  [1, 4, \"12+22=\", \"15+68=\"],\n  [2, 0, \"90-70=\", \"23+25=\"],\n  [2, 1, \"72-24=\", \"53-19=\"],\n  [2, 2, \"90-75=\", \"56+30=\"],\n  [2, 3, \"63-50=\", \"41+48=\"],\n  [2, 4, \"31-4=\", \"85-55=\"],\n  [3, 0, \"99-38=\", \"67+2=\"],\n  [3, 1, \"73+3=\", \"49-21=\"],\n  [3, 2, \"48+46=\", \"31+54=\"],\n  [3, 3, \"86-43=\", \"44+29=\"],\n  [3, 4, \"78-8=\", \"42-8=\"],\n  [4, 0, \"16+72=\", \"8+65=\"],\n  [4, 1, \"88-47=\", \"86-19=\"],\n  [4, 2, \"39+10=\", \"4+1=\"],\n  [4, 3, \"8+36=\", \"38+9=\"],\n  [4, 4, \"36-0=\", \"13+37=\"],\n  [5, 0, \"57+20=\", \"42-17=\"],\n  [5, 1, \"39+59=\", \"6+18=\"],\n  [5, 2, \"42+36=\", \"59+17=\"],\n  [5, 3, \"77-41=\", \"30-29=\"],\n  [5, 4, \"35-3=\", \"0+46=\"],\n  [6, 0, \"0+26=\", \"78-20=\"],\n  [6, 1, \"57-43=\", \"34+18=\"],\n  [6, 2, \"33+41=\", \"83-22=\"],\n  [6, 3, \"75-58=\", \"83-39=\"],\n  [6, 4, \"26-12=\", \"95-8=\"],\n  [7, 0, \"36-21=\", \"53+7=\"],\n  [7, 1, \"23-11=\", \"83-45=\"],\n  [7, 2, \"7+23=\", \"86+4=\"],\n  [7, 3, \"23+73=\", \"11+36=\"],\n  [7, 4, \"93-86=\", \"22+41=\"],\n  [8, 0, \"84-1=\", \"32+59=\"],\n  [8, 1, \"55-52=\", \"5+84=\"],\n  [8, 2, \"65-7=\", \"56+30=\"],\n  [8, 3, \"36-26=\", \"70-17=\"],\n  [8, 4, \"99-18=\", \"29+67=\"],\n  [9, 0, \"27+32=\", \"8+86=\"],\n  [9, 1, \"65-7=\", \"58-16=\"],\n  [9, 2, \"61-28=\", \"26+23=\"],\n  [9, 3, \"90-36=\", \"49-14=\"],\n  [9, 4, \"87-73=\", \"55-27=\"],\n  [10, 0, \"74-27=\", \"21-10=\"],\n  [10, 1, \"41-18=\", \"43-0=\"],\n  [10, 2, \"37+20=\", \"31-15=\"],\n  [10, 3, \"52-16=\", \"8+56=\"],\n  [10, 4, \"2+6=\", \"4+71=\"],\n  [11, 0, \"6+39=\", \"93-8=\"],\n  [11, 1, \"4+18=\", \"76-31=\"],\n  [11, 2, \"82-69=\", \"71-50=\"],\n  [11, 3, \"1+80=\", \"12+73=\"],\n  [11, 4, \"41-3=\", \"45-28=\"],\n  [12, 0, \"9+2=\", \"54-28=\"],\n  [12, 1, \"27-14=\", \"58-0=\"],\n  [12, 2, \"59+18=\", \"24+5=\"],\n  [12, 3, \"48+27=\", \"29-13=\"],\n  [12, 4, \"72-55=\", \"81-59=\"],\n  [13, 0, \"55+16=\", \"35+8=\"],\n  [13, 1, \"73-36=\", \"14-11=\"],\n  [13, 2, \"39+24=\", \"14+23=\"],\n  [13, 3, \"52+17=\", \"77-64=\"],\n  [13, 4, \"22-3=\", \"30+44=\"],\n  [14, 0, \"68+5=\", \"74-32=\"],\n  [14, 1, \"57+17=\", \"54-49=\"],\n  [14, 2, \"52-40=\", \"92-34=\"],\n  [14, 3, \"63+35=\", \"57-54=\"],\n  [14, 4, \"6+5=\", \"86-44=\"],\n  [15, 0, \"61+30=\", \"11-0=\"],\n  [15, 1, \"87-82=\", \"0+93=\"],\n  [15, 2, \"41+6=\", \"35-12=\"],\n  [15, 3, \"77-56=\", \"72-1=\"],\n  [15, 4, \"94-62=\", \"0+49=\"],\n  [16, 0, \"74-40=\", \"91-82=\"],\n  [16, 1, \"36-5=\", \"57+4=\"],\n  [16, 2, \"95-24=\", \"87+10=\"],\n  [16, 3, \"2+56=\", \"59+29=\"],\n  [16, 4, \"46+45=\", \"6+37=\"],\n  [17, 0, \"3+82=\", \"66+25=\"],\n  [17, 1, \"5+44=\", \"72-51=\"],\n  [17, 2, \"17+49=\", \"14+6=\"],\n  [17, 3, \"28-0=\", \"92-68=\"],\n  [17, 4, \"5+85=\", \"7+74=\"],\n  [18, 0, \"7+85=\", \"20-8=\"],\n  [18, 1, \"94-18=\", \"35-19=\"],\n  [18, 2, \"83-62=\", \"10+57=\"],\n  [18, 3, \"66+30=\", \"59+36=\"],\n  [18, 4, \"39+18=\", \"93-76=\"],\n  [19, 0, \"81-48=\", \"7+91=\"],\n  [19, 1, \"2+38=\", \"47+49=\"],\n  [19, 2, \"16+73=\", \"45+22=\"],\n  [19, 4, \"25+12=\", \"27-12=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\nconst table = tables.items[0];\n\nfor (const [row, col, oldText, newText] of edits) {\n  const cell = table.getCell(row, col);\n  const searchResults = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length === 0) {\n    throw new Error(`Could not find \"${oldText}\" in cell (${row}, ${col})`);\n  }\n  // Replace just the matched text in-place so the run's formatting is kept.\n  searchResults.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the math expressions in the worksheet table, cell by cell,\n# addressing each cell by its (Row, Col) position (1-based, as COM expects)\n# so duplicate expression text (e.g. \"65-7=\" appearing twice) is handled\n# correctly. Setting Cell.Range.Text in place preserves the existing run\n# formatting (font, size) and paragraph properties (justification).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$edits = @(\n    @{ Row = 1; Col = 1; Old = \"75+19=\"; New = \"66+25=\" },\n    @{ Row = 1; Col = 2; Old = \"99-43=\"; New = \"72-6=\" },\n    @{ Row = 1; Col = 3; Old = \"89-3=\"; New = \"2+97=\" },\n    @{ Row = 1; Col = 4; Old = \"12+66=\"; New = \"79+17=\" },\n    @{ Row = 1; Col = 5; Old = \"90-6=\"; New = \"7+67=\" },\n    @{ Row = 2; Col = 1; Old = \"4+25=\"; New = \"93-24=\" },\n    @{ Row = 2; Col = 2; Old = \"89-78=\"; New = \"39+29=\" },\n    @{ Row = 2; Col = 3; Old = \"74-25=\"; New = \"82-1=\" },\n    @{ Row = 2; Col = 4; Old = \"53-34=\"; New = \"46+51=\" },\n    @{ Row = 2; Col = 5; Old = \"12+22=\"; New = \"15+68=\" },\n    @{ Row = 3; Col = 1; Old = \"90-70=\"; New = \"23+25=\" },\n    @{ Row = 3; Col = 2; Old = \"72-24=\"; New = \"53-19=\" },\n    @{ Row = 3; Col = 3; Old = \"90-75=\"; New = \"56+30=\" },\n    @{ Row = 3; Col = 4; Old = \"63-50=\"; New = \"41+48=\" },\n    @{ Row = 3; Col = 5; Old = \"31-4=\"; New = \"85-55=\" },\n    @{ Row = 4; Col = 1; Old = \"99-38=\"; New = \"67+2=\" },\n    @{ Row = 4; Col = 2; Old = \"73+3=\"; New = \"49-21=\" },\n    @{ Row = 4; Col = 3; Old = \"48+46=\"; New = \"31+54=\" },\n    @{ Row = 4; Col = 4; Old = \"86-43=\"; New = \"44+29=\" },\n    @{ Row = 4; Col = 5; Old = \"78-8=\"; New = \"42-8=\" },\n    @{ Row = 5; Col = 1; Old = \"16+72=\"; New = \"8+65=\" },\n    @{ Row = 5; Col = 2; Old = \"88-47=\"; New = \"86-19=\" },\n    @{ Row = 5; Col = 3; Old = \"39+10=\"; New = \"4+1=\" },\n    @{ Row = 5; Col = 4; Old = \"8+36=\"; New = \"38+9=\" },\n    @{ Row = 5; Col = 5; Old = \"36-0=\"; New = \"13+37=\" },\n    @{ Row = 6; Col = 1; Old = \"57+20=\"; New = \"42-17=\" },\n    @{ Row = 6; Col = 2; Old = \"39+59=\"; New = \"6+18=\" },\n    @{ Row = 6; Col = 3; Old = \"42+36=\"; New = \"59+17=\" },\n    @{ Row = 6; Col = 4; Old = \"77-41=\"; New = \"30-29=\" },\n    @{ Row = 6; Col = 5; Old = \"35-3=\"; New = \"0+46=\" },\n    @{ Row = 7; Col = 1; Old = \"0+26=\"; New = \"78-20=\" },\n    @{ Row = 7; Col = 2; Old = \"57-43=\"; New = \"34+18=\" },\n    @{ Row = 7; Col = 3; Old = \"33+41=\"; New = \"83-22=\" },\n    @{ Row = 7; Col = 4; Old = \"75-58=\"; New = \"83-39=\" },\n    @{ Row = 7; Col = 5; Old = \"26-12=\"; New = \"95-8=\" },\n    @{ Row = 8; Col = 1; Old = \"36-21=\"; New = \"53+7=\" },\n    @{ Row = 8; Col = 2; Old = \"23-11=\"; New = \"83-45=\" },\n    @{ Row = 8; Col = 3; Old = \"7+23=\"; New = \"86+4=\" },\n    @{ Row = 8; Col = 4; Old = \"23+73=\"; New = \"11+36=\" },\n    @{ Row = 8; Col = 5; Old = \"93-86=\"; New = \"22+41=\" },\n    @{ Row = 9; Col = 1; Old = \"84-1=\"; New = \"32+59=\" },\n    @{ Row = 9; Col = 2; Old = \"55-52=\"; New = \"5+84=\" },\n    @{ Row = 9; Col = 3; Old = \"65-7=\"; New = \"56+30=\" },\n    @{ Row = 9; Col = 4; Old = \"36-26=\"; New = \"70-17=\" },\n    @{ Row = 9; Col = 5; Old = \"99-18=\"; New = \"29+67=\" },\n    @{ Row = 10; Col = 1; Old = \"27+32=\"; New = \"8+86=\" },\n    @{ Row = 10; Col = 2; Old = \"65-7=\"; New = \"58-16=\" },\n    @{ Row = 10; Col = 3; Old = \"61-28=\"; New = \"26+23=\" },\n    @{ Row = 10; Col = 4; Old = \"90-36=\"; New = \"49-14=\" },\n    @{ Row = 10; Col = 5; Old = \"87-73=\"; New = \"55-27=\" },\n    @{ Row = 11; Col = 1; Old = \"74-27=\"; New = \"21-10=\" },\n    @{ Row = 11; Col = 2; Old = \"41-18=\"; New = \"43-0=\" },\n    @{ Row = 11; Col = 3; Old = \"37+20=\"; New = \"31-15=\" },\n    @{ Row = 11; Col = 4; Old = \"52-16=\"; New = \"8+56=\" },\n    @{ Row = 11; Col = 5; Old = \"2+6=\"; New = \"4+71=\" },\n    @{ Row = 12; Col = 1; Old = \"6+39=\"; New = \"93-8=\" },\n    @{ Row = 12; Col = 2; Old = \"4+18=\"; New = \"76-31=\" },\n    @{ Row = 12; Col = 3; Old = \"82-69=\"; New = \"71-50=\" },\n    @{ Row = 12; Col = 4; Old = \"1+80=\"; New = \"12+73=\" },\n    @{ Row = 12; Col = 5; Old = \"41-3=\"; New = \"45-28=\" },\n    @{ Row = 13; Col = 1; Old = \"9+2=\"; New = \"54-28=\" },\n    @{ Row = 13; Col = 2; Old = \"27-14=\"; New = \"58-0=\" },\n    @{ Row = 13; Col = 3; Old = \"59+18=\"; New = \"24+5=\" },\n    @{ Row = 13; Col = 4; Old = \"48+27=\"; New = \"29-13=\" },\n    @{ Row = 13; Col = 5; Old = \"72-55=\"; New = \"81-59=\" },\n    @{ Row = 14; Col = 1; Old = \"55+16=\"; New = \"35+8=\" },\n    @{ Row = 14; Col = 2; Old = \"73-36=\"; New = \"14-11=\" },\n    @{ Row = 14; Col = 3; Old = \"39+24=\"; New = \"14+23=\" },\n    @{ Row = 14; Col = 4; Old = \"52+17=\"; New = \"77-64=\" },\n    @{ Row = 14; Col = 5; Old = \"22-3=\"; New = \"30+44=\" },\n    @{ Row = 15; Col = 1; Old = \"68+5=\"; New = \"74-32=\" },\n    @{ Row = 15; Col = 2; Old = \"57+17=\"; New = \"54-49=\" },\n    @{ Row = 15; Col = 3; Old = \"52-40=\"; New = \"92-34=\" },\n    @{ Row = 15; Col = 4; Old = \"63+35=\"; New = \"57-54=\" },\n    @{ Row = 15; Col = 5; Old = \"6+5=\"; New = \"86-44=\" },\n    @{ Row = 16; Col = 1; Old = \"61+30=\"; New = \"11-0=\" },\n    @{ Row = 16; Col = 2; Old = \"87-82=\"; New = \"0+93=\" },\n    @{ Row = 16; Col = 3; Old = \"41+6=\"; New = \"35-12=\" },\n    @{ Row = 16; Col = 4; Old = \"77-56=\"; New = \"72-1=\" },\n    @{ Row = 16; Col = 5; Old = \"94-62=\"; New = \"0+49=\" },\n    @{ Row = 17; Col = 1; Old = \"74-40=\"; New = \"91-82=\" },\n    @{ Row = 17; Col = 2; Old = \"36-5=\"; New = \"57+4=\" },\n    @{ Row = 17; Col = 3; Old = \"95-24=\"; New = \"87+10=\" },\n    @{ Row = 17; Col = 4; Old = \"2+56=\"; New = \"59+29=\" },\n    @{ Row = 17; Col = 5; Old = \"46+45=\"; New = \"6+37=\" },\n    @{ Row = 18; Col = 1; Old = \"3+82=\"; New = \"66+25=\" },\n    @{ Row = 18; Col = 2; Old = \"5+44=\"; New = \"72-51=\" },\n    @{ Row = 18; Col = 3; Old = \"17+49=\"; New = \"14+6=\" },\n    @{ Row = 18; Col = 4; Old = \"28-0=\"; New = \"92-68=\" },\n    @{ Row = 18; Col = 5; Old = \"5+85=\"; New = \"7+74=\" },\n    @{ Row = 19; Col = 1; Old = \"7+85=\"; New = \"20-8=\" },\n    @{ Row = 19; Col = 2; Old = \"94-18=\"; New = \"35-19=\" },\n    @{ Row = 19; Col = 3; Old = \"83-62=\"; New = \"10+57=\" },\n    @{ Row = 19; Col = 4; Old = \"66+30=\"; New = \"59+36=\" },\n    @{ Row = 19; Col = 5; Old = \"39+18=\"; New = \"93-76=\" },\n    @{ Row = 20; Col = 1; Old = \"81-48=\"; New = \"7+91=\" },\n    @{ Row = 20; Col = 2; Old = \"2+38=\"; New = \"47+49=\" },\n    @{ Row = 20; Col = 3; Old = \"16+73=\"; New = \"45+22=\" },\n    @{ Row = 20; Col = 5; Old = \"25+12=\"; New = \"27-12=\" }\n)\n\n$updated = 0\nforeach ($edit in $edits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $edit.Old) {\n        Write-Output \"Warning: Cell($($edit.Row),$($edit.Col)) expected '$($edit.Old)' but found '$current'\"\n    }\n    $cell.Range.Text = $edit.New\n    $updated++\n}\n\nWrite-Output \"Updated $updated cells\"\n"}
